$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.034.96"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.649.62"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.29"
$ws.Range("E5").Value = "  +3.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5233"
$ws.Range("E6").Value = "  +1.51%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2620"
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06361"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.74"
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07731"
$ws.Range("E11").Value = "  +2.72%  "
$ws.Range("D12").Value = "1.643.03"
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.444"
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("D14").Value = "1.874.03"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5502"
$ws.Range("E15").Value = "  +2.34%  "
$ws.Range("E16").Value = "  +4.50%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.79"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "26.042.83"
$ws.Range("E18").Value = "  -0.23%  "
$ws.Range("E19").Value = "  -0.19%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.737"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.60"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.334"
$ws.Range("E23").Value = "  +2.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.53"
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1247"
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.389"
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.01"
$ws.Range("E28").Value = "  +2.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.419"
$ws.Range("E29").Value = "  +2.81%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05938"
$ws.Range("E30").Value = "  -3.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.260"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.429"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.408"
$ws.Range("E33").Value = "  +0.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.650"
$ws.Range("E34").Value = "  +1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9924"
$ws.Range("E35").Value = "  +0.82%  "
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.757"
$ws.Range("E37").Value = "  +0.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5624"
$ws.Range("E38").Value = "  -4.22%  "
$ws.Range("E39").Value = "  +0.89%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.872"
$ws.Range("E40").Value = "  -1.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8574"
$ws.Range("E41").Value = "  +1.47%  "
$ws.Range("E42").Value = "  -0.15%  "
$ws.Range("D43").Value = "1.026.27"
$ws.Range("E43").Value = "  -7.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.08"
$ws.Range("E44").Value = "  -0.86%  "
$ws.Range("D45").Value = "1.797.21"
$ws.Range("E45").Value = "  -0.80%  "
$ws.Range("D46").Value = "0.0₈107"
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.74"
$ws.Range("E47").Value = "  +1.71%  "
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.070"
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05148"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4208"
$ws.Range("E51").Value = "  -0.82%  "
